# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 258 in the "Zapallo italiano" /
# Macroferia Regional de Talca sheet, pushing the existing rows 258-312
# down to 259-313 (dimension grows from A1:R312 to A1:R313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 258; everything below shifts down by one.
$ws.Rows(258).Insert()

# Populate the new row with the new weekly record.
$ws.Range("A258").Value = 5
$ws.Range("B258").Value = "Macroferia Regional de Talca"
$ws.Range("C258").Value = "Maule"
$ws.Range("D258").Value = 44641
$ws.Range("E258").Value = 7
$ws.Range("F258").Value = 100112032
$ws.Range("G258").Value = "Zapallo italiano"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 500
$ws.Range("K258").Value = 6000
$ws.Range("L258").Value = 6000
$ws.Range("M258").Value = 6000
$ws.Range("N258").Value = "$/caja 50 unidades"
$ws.Range("O258").Value = "Región del Maule"
$ws.Range("P258").Value = 120
$ws.Range("Q258").Value = 50
$ws.Range("R258").Value = "Hortaliza"

# Match the style used by the other date cells in column D (custom
# date/time number format), since Rows.Insert already copies the
# formatting of the row above but we make it explicit just in case.
$ws.Range("D258").NumberFormat = $ws.Range("D259").NumberFormat
